$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1673.625
$ws.Range("I19").Value = 1574.5
$ws.Range("J19").Value = 1772.75
$ws.Range("K19").Value = 1574.5
$ws.Range("L19").Value = 1772.75
$ws.Range("M19").Value = -1399.5
$ws.Range("N19").Value = -2122.75
$ws.Range("H111").Value = 687
$ws.Range("J111").Value = 791
$ws.Range("L111").Value = 2373
$ws.Range("N111").Value = -8507
$ws.Range("H137").Value = 25118.91
$ws.Range("I137").Value = 29604.676
$ws.Range("K137").Value = 88814.02799999999
$ws.Range("M137").Value = -86264.02799999999
$ws.Range("H138").Value = 3054.7556
$ws.Range("I138").Value = 2955
$ws.Range("J138").Value = 3150.1738
$ws.Range("K138").Value = 8865
$ws.Range("L138").Value = 9450.5214
$ws.Range("M138").Value = -3725
$ws.Range("N138").Value = -19730.5214

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20850.424
$ws.Range("I32").Value = 22445.4
$ws.Range("K32").Value = 22445.4
$ws.Range("M32").Value = -22158.4
$ws.Range("H122").Value = 2135.7273
$ws.Range("I122").Value = 2135.7273
$ws.Range("K122").Value = 6407.1819
$ws.Range("M122").Value = -3957.1819
$ws.Range("H132").Value = 29625.63
$ws.Range("I132").Value = 31891
$ws.Range("J132").Value = 3196.3333
$ws.Range("K132").Value = 95673
$ws.Range("L132").Value = 9588.999899999999
$ws.Range("M132").Value = -93143
$ws.Range("N132").Value = -14648.9999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4127.706
$ws.Range("I86").Value = 4141.4165
$ws.Range("J86").Value = 4094.8
$ws.Range("K86").Value = 4141.4165
$ws.Range("L86").Value = 4094.8
$ws.Range("M86").Value = -3018.4165
$ws.Range("N86").Value = -6340.8
$ws.Range("H89").Value = 4127.706
$ws.Range("I89").Value = 4141.4165
$ws.Range("J89").Value = 4094.8
$ws.Range("K89").Value = 20707.0825
$ws.Range("L89").Value = 20474
$ws.Range("M89").Value = -15091.0825
$ws.Range("N89").Value = -31706
$ws.Range("H105").Value = 115061.555
$ws.Range("I105").Value = 4034.9167
$ws.Range("K105").Value = 4034.9167
$ws.Range("M105").Value = -2287.9167
$ws.Range("H134").Value = 2604.585
$ws.Range("I134").Value = 2463.1633
$ws.Range("J134").Value = 4337
$ws.Range("K134").Value = 7389.4899
$ws.Range("L134").Value = 13011
$ws.Range("M134").Value = -4854.4899
$ws.Range("N134").Value = -18081

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H22").Value = 929
$ws.Range("J22").Value = 919.5714
$ws.Range("L22").Value = 919.5714
$ws.Range("N22").Value = -1619.5714
$ws.Range("H31").Value = 2273
$ws.Range("I31").Value = 1470.0667
$ws.Range("K31").Value = 1470.0667
$ws.Range("M31").Value = -1175.0667
$ws.Range("H34").Value = 2273
$ws.Range("I34").Value = 1470.0667
$ws.Range("K34").Value = 1470.0667
$ws.Range("M34").Value = -1268.0667
$ws.Range("H99").Value = 1989.5
$ws.Range("I99").Value = 1800
$ws.Range("K99").Value = 1800
$ws.Range("M99").Value = -302
$ws.Range("H107").Value = 445.25
$ws.Range("I107").Value = 490
$ws.Range("K107").Value = 490
$ws.Range("M107").Value = 1430
$ws.Range("H126").Value = 1989.5
$ws.Range("I126").Value = 1800
$ws.Range("K126").Value = 5400
$ws.Range("M126").Value = -2930
$ws.Range("H132").Value = 1631.9474
$ws.Range("I132").Value = 1631.9474
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4895.8422
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2365.8422
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 112775.11
$ws.Range("I134").Value = 112775.11
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 338325.33
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -335790.33
$ws.Range("N134").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 495
$ws.Range("I36").Value = 495
$ws.Range("K36").Value = 1485
$ws.Range("M36").Value = -1316
$ws.Range("H37").Value = 64950
$ws.Range("J37").Value = 64950
$ws.Range("L37").Value = 194850
$ws.Range("N37").Value = -195074
$ws.Range("H50").Value = 111323.664
$ws.Range("I50").Value = 171.33333
$ws.Range("K50").Value = 513.99999
$ws.Range("M50").Value = -32.99999000000003
$ws.Range("H53").Value = 111323.664
$ws.Range("I53").Value = 171.33333
$ws.Range("K53").Value = 513.99999
$ws.Range("M53").Value = -32.99999000000003
$ws.Range("H81").Value = 3500
$ws.Range("I81").Value = 3500
$ws.Range("K81").Value = 10500
$ws.Range("M81").Value = -9377
$ws.Range("H84").Value = 3500
$ws.Range("I84").Value = 3500
$ws.Range("K84").Value = 31500
$ws.Range("M84").Value = -25884
$ws.Range("H88").Value = 8586.77
$ws.Range("J88").Value = 9871.666999999999
$ws.Range("L88").Value = 29615.001
$ws.Range("N88").Value = -30471.001
$ws.Range("H91").Value = 8586.77
$ws.Range("J91").Value = 9871.666999999999
$ws.Range("L91").Value = 29615.001
$ws.Range("N91").Value = -32579.001
$ws.Range("H107").Value = 612.625
$ws.Range("J107").Value = 1200.6666
$ws.Range("L107").Value = 3601.9998
$ws.Range("N107").Value = -7441.9998
$ws.Range("H117").Value = 2080.2104
$ws.Range("I117").Value = 1002.5
$ws.Range("J117").Value = 2577.6155
$ws.Range("K117").Value = 3007.5
$ws.Range("L117").Value = 7732.8465
$ws.Range("M117").Value = 434.5
$ws.Range("N117").Value = -14616.8465
$ws.Range("H121").Value = 431.2
$ws.Range("I121").Value = 385
$ws.Range("K121").Value = 1155
$ws.Range("M121").Value = 155
$ws.Range("H131").Value = 4778159
$ws.Range("I131").Value = 15369.857
$ws.Range("J131").Value = 7159553.5
$ws.Range("K131").Value = 46109.571
$ws.Range("L131").Value = 21478660.5
$ws.Range("M131").Value = -41069.571
$ws.Range("N131").Value = -21488740.5
$ws.Range("H132").Value = 1484.3846
$ws.Range("I132").Value = 1341.4166
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 12072.7494
$ws.Range("L132").Value = 28800
$ws.Range("M132").Value = -9542.749400000001
$ws.Range("N132").Value = -33860
$ws.Range("H137").Value = 3614.0833
$ws.Range("I137").Value = 2096.6667
$ws.Range("J137").Value = 8166.3335
$ws.Range("K137").Value = 6290.000100000001
$ws.Range("L137").Value = 24499.0005
$ws.Range("M137").Value = -1190.000100000001
$ws.Range("N137").Value = -34699.00049999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6699.2856
$ws.Range("I70").Value = 6665.6665
$ws.Range("J70").Value = 6724.5
$ws.Range("K70").Value = 6665.6665
$ws.Range("L70").Value = 6724.5
$ws.Range("M70").Value = -6395.6665
$ws.Range("N70").Value = -7264.5
$ws.Range("H73").Value = 6699.2856
$ws.Range("I73").Value = 6665.6665
$ws.Range("J73").Value = 6724.5
$ws.Range("K73").Value = 6665.6665
$ws.Range("L73").Value = 6724.5
$ws.Range("M73").Value = -5729.6665
$ws.Range("N73").Value = -8596.5
$ws.Range("H100").Value = 34000
$ws.Range("J100").Value = 34000
$ws.Range("L100").Value = 34000
$ws.Range("N100").Value = -36164
$ws.Range("H122").Value = 5732.3335
$ws.Range("I122").Value = 3598.75
$ws.Range("K122").Value = 10796.25
$ws.Range("M122").Value = -8346.25
$ws.Range("H126").Value = 7752.1665
$ws.Range("I126").Value = 7151.8184
$ws.Range("K126").Value = 21455.4552
$ws.Range("M126").Value = -18985.4552
$ws.Range("H132").Value = 21652.963
$ws.Range("I132").Value = 32139.146
$ws.Range("K132").Value = 96417.43799999999
$ws.Range("M132").Value = -93887.43799999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7345.727
$ws.Range("I40").Value = 5127.2666
$ws.Range("K40").Value = 5127.2666
$ws.Range("M40").Value = -4991.2666
$ws.Range("H132").Value = 57740.727
$ws.Range("I132").Value = 69372.89
$ws.Range("K132").Value = 208118.67
$ws.Range("M132").Value = -205588.67
$ws.Range("H136").Value = 2254.516
$ws.Range("I136").Value = 2105.3809
$ws.Range("K136").Value = 6316.1427
$ws.Range("M136").Value = -3766.1427

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 11333.333
$ws.Range("I4").Value = 16000
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -15887
$ws.Range("N4").Value = -2226
$ws.Range("H5").Value = 5296
$ws.Range("I5").Value = 6000
$ws.Range("J5").Value = 4944
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 4944
$ws.Range("M5").Value = -5888
$ws.Range("N5").Value = -5168
$ws.Range("H132").Value = 22231.412
$ws.Range("I132").Value = 22231.412
$ws.Range("K132").Value = 66694.236
$ws.Range("M132").Value = -64164.236

Write-Host "Applied all profit sheet updates"